# doc: add instructions for input
#
# Rename the default sheet to "fixed", add two header/instruction rows
# ("_id" / "auto_increment"), size column A to fit the new content, and
# move the active selection down to the next empty row (A3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet1 -> fixed
$ws.Name = "fixed"

# Populate the instructional rows (these become xl/sharedStrings.xml
# entries "_id" / "auto_increment" automatically).
$ws.Range("A1").Value = "_id"
$ws.Range("A2").Value = "auto_increment"

# Widen column A so the longer string ("auto_increment") fits.
$ws.Columns.Item(1).ColumnWidth = 13.330729166666666

# Leave the cursor on the next free row, like a user who just finished
# typing the two instruction cells.
[void]$ws.Range("A3").Select()
